$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.361.29'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.533.95'
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("D4").Formula = "'0.998"
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Formula = "'544.01"
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D6").Formula = "'144.94"
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("D7").Formula = "'0.995"
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").Formula = "'0.573"
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("D9").Value = '2.553.25'
$ws.Range("E9").Value = '  +3.61%  '
$ws.Range("E10").Value = '  +2.09%  '
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Formula = "'5.60"
$ws.Range("E12").Value = '  +4.98%  '
$ws.Range("E13").Value = '  +1.30%  '
$ws.Range("D14").Value = '2.972.03'
$ws.Range("E14").Value = '  +2.96%  '
$ws.Range("D15").Formula = "'23.86"
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '59.319.33'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Formula = "'0.0000142"
$ws.Range("E17").Value = '  +2.79%  '
$ws.Range("D18").Value = '2.552.09'
$ws.Range("E18").Value = '  +1.41%  '
$ws.Range("E19").Value = '  +1.71%  '
$ws.Range("E20").Value = '  -0.85%  '
$ws.Range("D21").Formula = "'327.42"
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("E22").Value = '  +3.12%  '
$ws.Range("E23").Value = '  +2.82%  '
$ws.Range("D24").Formula = "'62.04"
$ws.Range("E24").Value = '  +2.20%  '
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("E26").Value = '  +2.51%  '
$ws.Range("D27").Formula = "'0.992"
$ws.Range("E27").Value = '  +1.57%  '
$ws.Range("D29").Formula = "'6.91"
$ws.Range("E29").Value = '  +3.28%  '
$ws.Range("D30").Value = '0.0₃0785'
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("E31").Value = '  +0.35%  '
$ws.Range("D32").Formula = "'1.22"
$ws.Range("E32").Value = '  -2.68%  '
$ws.Range("E33").Value = '  +9.83%  '
$ws.Range("D34").Formula = "'0.996"
$ws.Range("E34").Value = '  -0.11%  '
$ws.Range("D35").Formula = "'156.65"
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  -1.05%  '
$ws.Range("D38").Formula = "'1.62"
$ws.Range("E38").Value = '  -4.97%  '
$ws.Range("D39").Formula = "'5.71"
$ws.Range("E39").Value = '  -2.11%  '
$ws.Range("D40").Formula = "'37.13"
$ws.Range("E40").Value = '  +2.57%  '
$ws.Range("D41").Formula = "'300.25"
$ws.Range("E41").Value = '  -5.30%  '
$ws.Range("D42").Formula = "'3.72"
$ws.Range("D43").Formula = "'0.833"
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("D44").Formula = "'0.995"
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Formula = "'0.607"
$ws.Range("E45").Value = '  +4.60%  '
$ws.Range("D46").Formula = "'10.77"
$ws.Range("E46").Value = '  +0.35%  '
$ws.Range("D47").Formula = "'0.0935"
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("E48").Value = '  +2.43%  '
$ws.Range("D49").Formula = "'124.11"
$ws.Range("E49").Value = '  +2.00%  '
$ws.Range("E51").Value = '  -1.90%  '
